# New crime data collected - weekly CompStat update
# Updates the report header (volume/week dates) and the weekly crime-stat
# table (rows 15-30) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# ---------------------------------------------------------------------
# Header: Volume/Number and the week-covering date range
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# ---------------------------------------------------------------------
# Helper-free, explicit per-cell updates for the crime stats table.
# Numeric formats (NumberFormat) are (re)applied only where the source
# cell used to hold a text placeholder ("n/a" / "***.*") and now needs
# to carry a real number, so the stored style switches from the
# text-style to the matching numeric style, exactly like Excel does
# when you type a number into a cell that previously held text.
# ---------------------------------------------------------------------

$fmtInt = "#,##0"
$fmtPct = "#,##0.0;`"-`"#,##0.0"

# Row 15 - Rape
$ws.Range("L15").Value = 14.285714285714
$ws.Range("N15").Value = -55.555555555555

# Row 16 - Robbery
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 46
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = -16.363636363636
$ws.Range("L16").Value = 27.777777777777
$ws.Range("M16").Value = -52.577319587628
$ws.Range("N16").Value = -83.150183150183

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 151
$ws.Range("J17").Value = 115
$ws.Range("K17").Value = 31.304347826087
$ws.Range("L17").Value = 36.036036036036
$ws.Range("M17").Value = 20.8
$ws.Range("N17").Value = -45.878136200716

# Row 18 - Burglary
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 82
$ws.Range("J18").Value = 71
$ws.Range("K18").Value = 15.492957746478
$ws.Range("L18").Value = 24.242424242424
$ws.Range("M18").Value = -56.149732620320
$ws.Range("N18").Value = -93.234323432343

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -13.333333333333
$ws.Range("I19").Value = 425
$ws.Range("J19").Value = 322
$ws.Range("K19").Value = 31.987577639751
$ws.Range("L19").Value = 44.557823129251
$ws.Range("M19").Value = 15.803814713896
$ws.Range("N19").Value = -42.876344086021

# Row 20 - G.L.A.
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 99
$ws.Range("J20").Value = 124
$ws.Range("K20").Value = -20.161290322580
$ws.Range("L20").Value = 39.436619718309
$ws.Range("M20").Value = -5.714285714285
$ws.Range("N20").Value = -96.024096385542

# Row 21 - TOTAL
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 21.428571428571
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -21.176470588235
$ws.Range("I21").Value = 813
$ws.Range("J21").Value = 696
$ws.Range("K21").Value = 16.810344827586
$ws.Range("L21").Value = 38.500851788756
$ws.Range("M21").Value = -9.566184649610
$ws.Range("N21").Value = -83.801554094441

# Row 23 - Housing
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 32
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = 68.421052631578
$ws.Range("L23").Value = -8.571428571428
$ws.Range("M23").Value = 52.380952380952

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 11.764705882352
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 79
$ws.Range("H24").Value = 11.392405063291
$ws.Range("I24").Value = 966
$ws.Range("J24").Value = 764
$ws.Range("K24").Value = 26.439790575916
$ws.Range("L24").Value = 117.078651685393
$ws.Range("M24").Value = -36.279683377308

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 31.034482758620
$ws.Range("I25").Value = 297
$ws.Range("J25").Value = 301
$ws.Range("K25").Value = -1.328903654485
$ws.Range("L25").Value = 35
$ws.Range("M25").Value = -38.125

# Row 26 - UCR Rape*  (D26/E26 flip from the "n/a"/"***.*" placeholder to
# real numbers, so their NumberFormat needs to be (re)applied to land on
# the correct numeric style.)
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = $fmtInt
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = $fmtPct
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -22.222222222222
$ws.Range("L26").Value = 0

# Row 27 - Other Sex Crimes (C27/D27/E27/G27/H27 flip from placeholder
# text to real numbers.)
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = $fmtInt
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = $fmtInt
$ws.Range("E27").Value = -66.666666666666
$ws.Range("E27").NumberFormat = $fmtPct
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("G27").NumberFormat = $fmtInt
$ws.Range("H27").Value = 66.666666666666
$ws.Range("H27").NumberFormat = $fmtPct
$ws.Range("I27").Value = 34
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = 17.241379310344
$ws.Range("L27").Value = 13.333333333333

# Row 30 - Hate Crimes (F30 flips from placeholder text to a real number.)
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = $fmtInt
$ws.Range("I30").Value = 8
$ws.Range("K30").Value = 14.285714285714
$ws.Range("L30").Value = 100
